# with price above and below operations
# Adds a new "Section 3" block (rows 26-37) describing operation/signal
# columns used for crossAbove / topLine / variable-crossPercent signals,
# widens columns D:F to fit the new content and scrolls/selects near the
# new block the way Excel would after a user finished typing it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New "Section 3" block --------------------------------------------
$ws.Range("C26").Value = "Section 3"

$ws.Range("C27").Value = "operation"
$ws.Range("D27").Value = "signal"
$ws.Range("E27").Value = "signal"
$ws.Range("F27").Value = "signal"

$ws.Range("B28").Value = "Col 1"
$ws.Range("D28").Value = "variable crossPercent 3"

$ws.Range("B29").Value = "Col 2"
$ws.Range("D29").Value = "crossAbove 10 50"
$ws.Range("E29").Value = "crossAbove 10 100"
$ws.Range("F29").Value = "crossAbove 10 200"

$ws.Range("B30").Value = "Col 3"
$ws.Range("D30").Value = "topLine 10"

$ws.Range("B31").Value = "Col 4"
$ws.Range("B32").Value = "Col 5"
$ws.Range("B33").Value = "Col 6"
$ws.Range("B34").Value = "Col 7"
$ws.Range("B35").Value = "Col 8"
$ws.Range("B36").Value = "Col 9"
$ws.Range("B37").Value = "Col 10"

# Match the row height used for the rest of the new block.
$ws.Range("B26:F37").RowHeight = 13.8

# ---- Column widths ------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.79
$ws.Columns.Item(5).ColumnWidth = 19.34
$ws.Columns.Item(6).ColumnWidth = 20.78

# ---- Selection / scroll position ----------------------------------------
$ws.Range("D33").Select()
